$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 8 (even_MAG-GUT84184.fa), shifting rows 9-10 up to 8-9
$ws.Rows.Item(8).Delete()
